$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44635
$ws.Range("J2").Value = 100

# Row 3
$ws.Range("D3").Value = 44658
$ws.Range("J3").Value = 80

# Row 4
$ws.Range("D4").Value = 44664
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 160
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 861

# Row 5
$ws.Range("D5").Value = 44651
$ws.Range("J5").Value = 60

# Row 6
$ws.Range("D6").Value = 44659
$ws.Range("J6").Value = 80

# Row 7
$ws.Range("D7").Value = 44637
$ws.Range("J7").Value = 100

# Row 9
$ws.Range("D9").Value = 44384
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 60
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 15000
$ws.Range("P9").Value = 833

# Row 10
$ws.Range("D10").Value = 44642
$ws.Range("J10").Value = 100

# Row 11
$ws.Range("D11").Value = 44628
$ws.Range("J11").Value = 60

# Row 12
$ws.Range("D12").Value = 44649

# Row 13
$ws.Range("D13").Value = 44645
$ws.Range("J13").Value = 60

# Row 14
$ws.Range("D14").Value = 44630
